$wb = $excel.ActiveWorkbook

# Sheets
$lcv = $wb.Worksheets.Item("LCV")
$binek = $wb.Worksheets.Item("binek_arac")

# Add new sheet "HDV" after LCV
$hdv = $wb.Worksheets.Add($null, $lcv)
$hdv.Name = "HDV"

# Fill header + data, mirroring LCV's layout
$hdv.Range("A1").Value = "degisken"
$hdv.Range("B1").Value = "deger"
$hdv.Range("A1:B1").Font.Bold = $true

$hdv.Range("A2").Value = "HDV_hurda_tesvik_orani"
$hdv.Range("B2").Value = 0.15

$hdv.Columns.Item(1).ColumnWidth = 21.8

# Selections to match target state
$binek.Range("E7").Select() | Out-Null
$lcv.Range("A1:B2").Select() | Out-Null
$hdv.Range("H27").Select() | Out-Null

$hdv.Activate() | Out-Null
